# Generate Report for Handback
# The "cbea7538-5dfa-42de-bae8-bc405818f8b7" file has now been handed back
# (both zh-cn and de-de targets). Update the per-language status sheets and
# the Overview roll-up sheet to reflect this.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"
$mdName     = "cbea7538-5dfa-42de-bae8-bc405818f8b7.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b15c19b21d77f9e52f7d46941b9a9c8cbaaa7098/e2e/cbea7538-5dfa-42de-bae8-bc405818f8b7.md"
$hyperlinkColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's existing HyperLink style

# ---------------------------------------------------------------------------
# zh-cn sheet: row 6 is the cbea7538... entry that just got handed back
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C6").Value = $handedBack

$xlfZh = $wsZh.Range("G6").Value2
$wsZh.Range("K6").Value = $xlfZh

$wsZh.Range("L6").Value = "2017-02-09 08:41:57"
$wsZh.Range("L6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("J6"), $mdUrl, $null, $null, $mdName) | Out-Null
$wsZh.Range("J6").Font.Underline = $true
$wsZh.Range("J6").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# de-de sheet: row 6 is the cbea7538... entry that just got handed back
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C6").Value = $handedBack

$xlfDe = $wsDe.Range("G6").Value2
$wsDe.Range("K6").Value = $xlfDe

$wsDe.Range("L6").Value = "2017-02-09 08:42:24"
$wsDe.Range("L6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("J6"), $mdUrl, $null, $null, $mdName) | Out-Null
$wsDe.Range("J6").Font.Underline = $true
$wsDe.Range("J6").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Overview sheet: row 6 (cbea7538...) status columns for each language
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = $handedBack
$wsOverview.Range("F6").Value = $handedBack
